$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper donors (row 14 is untouched by this edit, so it is a stable source
# of already-existing cell styles/types that we need to re-apply elsewhere):
#   $donorText0    -> style 14 (General, right/center aligned), text "0"
#   $donorTextStar -> style 14 (General, right/center aligned), text "***.*"
#   $donorNum      -> style 15 (#,##0), number
#   $donorPct      -> style 16 (#,##0.0;"-"#,##0.0), number
# Copying a donor cell onto a target cell (Range.Copy(Destination)) clones
# both its value AND its style without ever minting a brand-new style record.
# ---------------------------------------------------------------------------
$donorText0    = $ws.Range("C14")
$donorTextStar = $ws.Range("E14")
$donorNum      = $ws.Range("F14")
$donorPct      = $ws.Range("M14")

# ----- Shared string text fixups (volume/issue number + report week dates) -
$ws.Range("A8").Value = "Volume 30   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/11/2023  Through  9/17/2023"

# ----- Row 15 (Murder) ------------------------------------------------------
$donorText0.Copy($ws.Range("C15"))
$ws.Range("M15").Value = 0

# ----- Row 16 (Rape) --------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 63
$ws.Range("J16").Value = 84
$ws.Range("K16").Value = -25
$ws.Range("L16").Value = -4.545454545454
$ws.Range("M16").Value = -5.970149253731
$ws.Range("N16").Value = -85.382830626450

# ----- Row 17 (Robbery) -----------------------------------------------------
$donorNum.Copy($ws.Range("C17"))
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -41.666666666666
$ws.Range("I17").Value = 79
$ws.Range("J17").Value = 76
$ws.Range("K17").Value = 3.947368421052
$ws.Range("L17").Value = 83.720930232558
$ws.Range("M17").Value = 68.085106382978
$ws.Range("N17").Value = 16.176470588235

# ----- Row 18 (Fel. Assault) ------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("I18").Value = 84
$ws.Range("J18").Value = 109
$ws.Range("K18").Value = -22.935779816513
$ws.Range("L18").Value = 42.372881355932
$ws.Range("M18").Value = 9.090909090909
$ws.Range("N18").Value = -86.895475819032

# ----- Row 19 (Burglary) ----------------------------------------------------
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -7.142857142857
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = -12.307692307692
$ws.Range("I19").Value = 519
$ws.Range("J19").Value = 586
$ws.Range("K19").Value = -11.433447098976
$ws.Range("L19").Value = 12.337662337662
$ws.Range("M19").Value = 11.373390557939
$ws.Range("N19").Value = -62.281976744186

# ----- Row 20 (Gr. Larceny) -------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -44.444444444444
$ws.Range("I20").Value = 74
$ws.Range("J20").Value = 56
$ws.Range("K20").Value = 32.142857142857
$ws.Range("L20").Value = 17.460317460317
$ws.Range("M20").Value = 236.363636363636
$ws.Range("N20").Value = -90.830235439900

# ----- Row 21 (G.L.A. / TOTAL) ----------------------------------------------
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = -14.851485148514
$ws.Range("I21").Value = 831
$ws.Range("J21").Value = 922
$ws.Range("K21").Value = -9.869848156182
$ws.Range("L21").Value = 18.884120171673
$ws.Range("M21").Value = 20.784883720930
$ws.Range("N21").Value = -75.134649910233

# ----- Row 22 (Transit) -----------------------------------------------------
$ws.Range("L22").Value = 0

# ----- Row 23 (Housing) -----------------------------------------------------
$ws.Range("I23").Value = 29
$ws.Range("K23").Value = 52.631578947368
$ws.Range("L23").Value = 81.25
$ws.Range("M23").Value = 38.095238095238

# ----- Row 24 (Petit Larceny) -----------------------------------------------
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 125
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = 32.978723404255
$ws.Range("I24").Value = 838
$ws.Range("J24").Value = 947
$ws.Range("K24").Value = -11.510031678986
$ws.Range("L24").Value = -11.134676564156
$ws.Range("M24").Value = 8.831168831168

# ----- Row 25 (Misd. Assault) -----------------------------------------------
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 46.153846153846
$ws.Range("I25").Value = 166
$ws.Range("J25").Value = 152
$ws.Range("K25").Value = 9.210526315789
$ws.Range("L25").Value = 43.103448275862
$ws.Range("M25").Value = -15.306122448979

# ----- Row 26 (UCR Rape*) ---------------------------------------------------
$donorText0.Copy($ws.Range("C26"))
$donorNum.Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$donorPct.Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 15
$ws.Range("K26").Value = 6.666666666666

# ----- Row 27 (Other Sex Crimes) --------------------------------------------
$donorText0.Copy($ws.Range("C27"))
$donorText0.Copy($ws.Range("D27"))
$donorTextStar.Copy($ws.Range("E27"))
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1

Write-Output "done"
